$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "71.167.73"
$ws.Range("E2").Value = "  +3.13%  "

# Row 3
$ws.Range("D3").Value = "3.816.97"
$ws.Range("E3").Value = "  +1.20%  "

# Row 4
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.03%  "

# Row 5
$ws.Range("D5").Value = "'707.45"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +11.85%  "

# Row 6
$ws.Range("D6").Value = "'173.67"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.96%  "

# Row 7
$ws.Range("D7").Value = "3.819.86"
$ws.Range("E7").Value = "  +1.29%  "

# Row 8
$ws.Range("E8").Value = "  +0.00%  "

# Row 9
$ws.Range("E9").Value = "  +1.32%  "

# Row 10
$ws.Range("E10").Value = "  +3.76%  "

# Row 11
$ws.Range("D11").Value = "'7.42"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +9.90%  "

# Row 12
$ws.Range("E12").Value = "  +1.54%  "

# Row 13
$ws.Range("D13").Value = "'0.0000258"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +8.21%  "

# Row 14
$ws.Range("D14").Value = "'36.36"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.44%  "

# Row 15
$ws.Range("D15").Value = "4.454.13"
$ws.Range("E15").Value = "  +1.10%  "

# Row 16
$ws.Range("D16").Value = "3.812.29"
$ws.Range("E16").Value = "  +1.04%  "

# Row 17
$ws.Range("D17").Value = "71.168.73"
$ws.Range("E17").Value = "  +3.16%  "

# Row 18
$ws.Range("D18").Value = "'18.07"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.38%  "

# Row 19
$ws.Range("D19").Value = "'7.27"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.69%  "

# Row 20
$ws.Range("E20").Value = "  +0.55%  "

# Row 21
$ws.Range("D21").Value = "'11.25"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +18.57%  "

# Row 22
$ws.Range("D22").Value = "'485.30"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.79%  "

# Row 23
$ws.Range("D23").Value = "'0.720"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.49%  "

# Row 24
$ws.Range("D24").Value = "'84.03"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.37%  "

# Row 25
$ws.Range("E25").Value = "  +2.42%  "

# Row 26
$ws.Range("D26").Value = "'12.56"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.97%  "

# Row 27
$ws.Range("D27").Value = "'10.67"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +5.50%  "

# Row 28
$ws.Range("D28").Value = "'2.20"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.97%  "

# Row 29
$ws.Range("D29").Value = "3.967.86"
$ws.Range("E29").Value = "  +1.04%  "

# Row 30
$ws.Range("E30").Value = "  -0.14%  "

# Row 31
$ws.Range("E31").Value = "  +15.22%  "

# Row 32
$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").Value = "'2.32"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.02%  "

# Row 33
$ws.Range("B33").Value = "NEARProtocol"
$ws.Range("C33").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D33").Value = "'7.61"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +7.55%  "

# Row 34
$ws.Range("D34").Value = "'29.74"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.65%  "

# Row 35
$ws.Range("D35").Value = "'0.180"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.06%  "

# Row 36
$ws.Range("E36").Value = "  +4.41%  "

# Row 37
$ws.Range("E37").Value = "  +0.05%  "

# Row 38
$ws.Range("D38").Value = "3.766.80"
$ws.Range("E38").Value = "  +1.16%  "

# Row 39
$ws.Range("E39").Value = "  +3.89%  "

# Row 40
$ws.Range("D40").Value = "'3.53"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +7.77%  "

# Row 41
$ws.Range("E41").Value = "  +4.29%  "

# Row 42
$ws.Range("B42").Value = "FLOKI"
$ws.Range("C42").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D42").Value = "'0.000346"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +31.42%  "

# Row 43
$ws.Range("B43").Value = "Stacks"
$ws.Range("C43").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D43").Value = "'2.26"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +14.15%  "

# Row 44
$ws.Range("D44").Value = "'0.974"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.34%  "

# Row 45
$ws.Range("E45").Value = "  -0.04%  "

# Row 46
$ws.Range("E46").Value = "  -0.01%  "

# Row 47
$ws.Range("D47").Value = "'45.71"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +5.07%  "

# Row 48
$ws.Range("B48").Value = "Monero"
$ws.Range("C48").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D48").Value = "'161.10"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.09%  "

# Row 49
$ws.Range("B49").Value = "OKB"
$ws.Range("C49").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D49").Value = "'49.45"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +5.36%  "

# Row 50
$ws.Range("E50").Value = "  +0.26%  "

# Row 51
$ws.Range("E51").Value = "  +3.04%  "
